$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 24 data for D,J,K,L,M,P)
$ws.Cells.Item(2,4).Value = (Get-Date -Year 2021 -Month 9 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(2,10).Value = 52
$ws.Cells.Item(2,11).Value = 5000
$ws.Cells.Item(2,12).Value = 6000
$ws.Cells.Item(2,13).Value = 5500
$ws.Cells.Item(2,16).Value = 344

# Row 3 (was row 5 data for D,J,K,L,M,P)
$ws.Cells.Item(3,4).Value = (Get-Date -Year 2021 -Month 4 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(3,10).Value = 70

# Row 4 (was row 23 data for D,J,K,L,M,P)
$ws.Cells.Item(4,4).Value = (Get-Date -Year 2021 -Month 6 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(4,10).Value = 43
$ws.Cells.Item(4,11).Value = 4500
$ws.Cells.Item(4,12).Value = 5000
$ws.Cells.Item(4,13).Value = 4756
$ws.Cells.Item(4,16).Value = 297

# Row 5 (was row 15 data for D,J,K,L,M,P)
$ws.Cells.Item(5,4).Value = (Get-Date -Year 2021 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(5,10).Value = 52
$ws.Cells.Item(5,11).Value = 5000
$ws.Cells.Item(5,13).Value = 5500
$ws.Cells.Item(5,16).Value = 344

# Row 6 (was row 7 data for D,J,K,L,M,P)
$ws.Cells.Item(6,4).Value = (Get-Date -Year 2021 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(6,10).Value = 50
$ws.Cells.Item(6,11).Value = 6000
$ws.Cells.Item(6,13).Value = 6000
$ws.Cells.Item(6,16).Value = 375

# Row 7 (was row 18 data for D,J,K,L,M,P)
$ws.Cells.Item(7,4).Value = (Get-Date -Year 2021 -Month 4 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(7,10).Value = 34

# Row 8 (was row 25 data for D,J,K,L,M,P)
$ws.Cells.Item(8,4).Value = (Get-Date -Year 2022 -Month 1 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(8,10).Value = 52
$ws.Cells.Item(8,11).Value = 8000
$ws.Cells.Item(8,12).Value = 8000
$ws.Cells.Item(8,13).Value = 8000
$ws.Cells.Item(8,16).Value = 500

# Row 9 (was row 16 data for D,J,K,L,M,P)
$ws.Cells.Item(9,4).Value = (Get-Date -Year 2021 -Month 10 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(9,10).Value = 25
$ws.Cells.Item(9,11).Value = 6000
$ws.Cells.Item(9,13).Value = 6000
$ws.Cells.Item(9,16).Value = 375

# Row 10 (was row 14 data for D,J,K,L,M,P)
$ws.Cells.Item(10,4).Value = (Get-Date -Year 2021 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(10,10).Value = 160
$ws.Cells.Item(10,12).Value = 6000
$ws.Cells.Item(10,13).Value = 6000
$ws.Cells.Item(10,16).Value = 375

# Row 11 (was row 3 data for D,J,K,L,M,P)
$ws.Cells.Item(11,4).Value = (Get-Date -Year 2021 -Month 6 -Day 3 -Hour 0 -Minute 0 -Second 0)

# Row 12 (was row 20 data for D,J,K,L,M,P)
$ws.Cells.Item(12,4).Value = (Get-Date -Year 2021 -Month 7 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(12,10).Value = 45
$ws.Cells.Item(12,13).Value = 5744
$ws.Cells.Item(12,16).Value = 359

# Row 13 (was row 9 data for D,J,K,L,M,P)
$ws.Cells.Item(13,4).Value = (Get-Date -Year 2021 -Month 6 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(13,10).Value = 160
$ws.Cells.Item(13,11).Value = 5500
$ws.Cells.Item(13,13).Value = 5750
$ws.Cells.Item(13,16).Value = 359

# Row 14 (was row 8 data for D,J,K,L,M,P)
$ws.Cells.Item(14,4).Value = (Get-Date -Year 2021 -Month 7 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(14,10).Value = 43

# Row 15 (was row 10 data for D,J,K,L,M,P)
$ws.Cells.Item(15,4).Value = (Get-Date -Year 2021 -Month 9 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(15,10).Value = 25
$ws.Cells.Item(15,11).Value = 6000
$ws.Cells.Item(15,12).Value = 7000
$ws.Cells.Item(15,13).Value = 6480
$ws.Cells.Item(15,16).Value = 405

# Row 16 (was row 12 data for D,J,K,L,M,P)
$ws.Cells.Item(16,4).Value = (Get-Date -Year 2021 -Month 5 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(16,10).Value = 51
$ws.Cells.Item(16,11).Value = 5500
$ws.Cells.Item(16,13).Value = 5755
$ws.Cells.Item(16,16).Value = 360

# Row 17 (was row 11 data for D,J,K,L,M,P)
$ws.Cells.Item(17,4).Value = (Get-Date -Year 2021 -Month 6 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(17,10).Value = 25
$ws.Cells.Item(17,11).Value = 6000
$ws.Cells.Item(17,12).Value = 6000
$ws.Cells.Item(17,13).Value = 6000
$ws.Cells.Item(17,16).Value = 375

# Row 18 (was row 17 data for D,J,K,L,M,P)
$ws.Cells.Item(18,4).Value = (Get-Date -Year 2022 -Month 1 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(18,10).Value = 61
$ws.Cells.Item(18,11).Value = 8000
$ws.Cells.Item(18,12).Value = 8000
$ws.Cells.Item(18,13).Value = 8000
$ws.Cells.Item(18,16).Value = 500

# Row 19 (was row 21 data for D,J,K,L,M,P)
$ws.Cells.Item(19,4).Value = (Get-Date -Year 2022 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(19,11).Value = 7000
$ws.Cells.Item(19,12).Value = 7000
$ws.Cells.Item(19,13).Value = 7000
$ws.Cells.Item(19,16).Value = 438

# Row 20 (was row 4 data for D,J,K,L,M,P)
$ws.Cells.Item(20,4).Value = (Get-Date -Year 2021 -Month 6 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(20,10).Value = 34
$ws.Cells.Item(20,13).Value = 5750

# Row 21 (was row 2 data for D,J,K,L,M,P)
$ws.Cells.Item(21,4).Value = (Get-Date -Year 2022 -Month 1 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(21,10).Value = 34
$ws.Cells.Item(21,11).Value = 8000
$ws.Cells.Item(21,12).Value = 8000
$ws.Cells.Item(21,13).Value = 8000
$ws.Cells.Item(21,16).Value = 500

# Row 22 (was row 13 data for D,J,K,L,M,P)
$ws.Cells.Item(22,4).Value = (Get-Date -Year 2021 -Month 5 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(22,10).Value = 120

# Row 23 (was row 6 data for D,J,K,L,M,P)
$ws.Cells.Item(23,4).Value = (Get-Date -Year 2021 -Month 8 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(23,10).Value = 34
$ws.Cells.Item(23,11).Value = 5000
$ws.Cells.Item(23,12).Value = 6000
$ws.Cells.Item(23,13).Value = 5500
$ws.Cells.Item(23,16).Value = 344

# Row 24 (was row 22 data for D,J,K,L,M,P)
$ws.Cells.Item(24,4).Value = (Get-Date -Year 2021 -Month 6 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(24,11).Value = 6000
$ws.Cells.Item(24,13).Value = 6000
$ws.Cells.Item(24,16).Value = 375

# Row 25 (was row 19 data for D,J,K,L,M,P)
$ws.Cells.Item(25,4).Value = (Get-Date -Year 2021 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(25,11).Value = 5000
$ws.Cells.Item(25,12).Value = 6000
$ws.Cells.Item(25,13).Value = 5500
$ws.Cells.Item(25,16).Value = 344
